$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MA and Decomp")
$ws.Range("K1").Value = "test"
